$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws2 = $wb.Worksheets.Item("tasas")

$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 3.04 = 11264.32 pesos`n✅ 11264.32 pesos = 3.02 = 958.23 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws2.Range("N10").Value = 329
$ws2.Range("O10").Value = 3705.96
$ws2.Range("N12").Value = 3729.95
$ws2.Range("O12").Value = 317.3
